$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel
# need an explicit Text number format so they remain strings, matching the
# original workbook convention of storing all Price/Volume entries as text.
$textCells = @('D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D13', 'D16', 'D19', 'D21', 'D22', 'D23', 'D26', 'D27', 'D30', 'D31', 'D32', 'D33', 'D35', 'D36', 'D38', 'D39', 'D42', 'D47', 'D48', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '42.831.20'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '2.542.65'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '303.78'
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('D6').Value = '97.78'
$ws.Range('E6').Value = '  +6.12%  '
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.545'
$ws.Range('E9').Value = '  -0.76%  '
$ws.Range('D10').Value = '36.84'
$ws.Range('E10').Value = '  +2.69%  '
$ws.Range('D11').Value = '0.0830'
$ws.Range('E11').Value = '  +3.32%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '7.72'
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.114'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = '2.934.64'
$ws.Range('E14').Value = '  +0.02%  '
$ws.Range('D15').Value = '2.536.75'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = '15.09'
$ws.Range('E16').Value = '  +5.98%  '
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').Value = '42.840.65'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('D19').Value = '13.34'
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('D21').Value = '6.58'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').Value = '71.80'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').Value = '255.98'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('E24').Value = '  +1.28%  '
$ws.Range('E25').Value = '  -1.62%  '
$ws.Range('D26').Value = '28.09'
$ws.Range('E26').Value = '  -3.60%  '
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('E28').Value = '  +9.26%  '
$ws.Range('E29').Value = '  +1.21%  '
$ws.Range('D30').Value = '37.93'
$ws.Range('E30').Value = '  +2.84%  '
$ws.Range('D31').Value = '6.08'
$ws.Range('E31').Value = '  +2.12%  '
$ws.Range('D32').Value = '157.71'
$ws.Range('E32').Value = '  +3.25%  '
$ws.Range('D33').Value = '19.47'
$ws.Range('E33').Value = '  +13.51%  '
$ws.Range('E34').Value = '  -1.48%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.0798'
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').Value = '3.31'
$ws.Range('E36').Value = '  -2.54%  '
$ws.Range('E37').Value = '  -4.24%  '
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = '25.90'
$ws.Range('E38').Value = '  +8.34%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '0.116'
$ws.Range('E39').Value = '  +1.42%  '
$ws.Range('E40').Value = '  +0.23%  '
$ws.Range('E41').Value = '  +29.34%  '
$ws.Range('D42').Value = '3.89'
$ws.Range('E42').Value = '  +0.59%  '
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('D44').Value = '2.091.36'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('E45').Value = '  -1.44%  '
$ws.Range('D47').Value = '87.68'
$ws.Range('E47').Value = '  +3.91%  '
$ws.Range('D48').Value = '8.86'
$ws.Range('E48').Value = '  -3.24%  '
$ws.Range('D49').Value = '2.792.34'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').Value = '74.49'
$ws.Range('E50').Value = '  +8.23%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.191'
$ws.Range('E51').Value = '  +2.71%  '
